{"js": "// Remove the \"Additional Notes:\" / \"{{ gig_notes or \"None\" }}\" block\n// (and its two trailing line breaks) from the performance-detail\n// paragraph, leaving \"Special Instructions:\" ... directly followed by\n// \"EVENT PACKAGE AND BAND SERVICES\".\nconst body = context.document.body;\n\nconst startResults = body.search(\"Additional Notes:\", { matchCase: true, matchWholeWord: false });\nconst endResults = body.search(\"EVENT PACKAGE AND BAND SERVICES\", { matchCase: true, matchWholeWord: false });\nstartResults.load(\"text\");\nendResults.load(\"text\");\nawait context.sync();\n\nif (startResults.items.length > 0 && endResults.items.length > 0) {\n  const startRange = startResults.items[0];\n  const endRange = endResults.items[0];\n\n  // Range running from the \"A\" in \"Additional Notes:\" up to (but not\n  // including) the start of \"EVENT PACKAGE AND BAND SERVICES\" \u2014 this\n  // covers \"Additional Notes:\" + the break + \"{{ gig_notes or \\\"None\\\" }}\"\n  // + the two trailing breaks.\n  const toDelete = startRange.expandTo(endRange.getRange(\"Before\"));\n  toDelete.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"Additional Notes:\" / \"{{ gig_notes or \"None\" }}\" block\n# (and its two trailing line breaks) from the performance-detail\n# paragraph, leaving \"Special Instructions:\" ... directly followed by\n# \"EVENT PACKAGE AND BAND SERVICES\".\n$d = $word.ActiveDocument\n\n$rngStart = $d.Content\n$foundStart = $rngStart.Find.Execute(\"Additional Notes:\")\n\n$rngEnd = $d.Content\n$foundEnd = $rngEnd.Find.Execute(\"EVENT PACKAGE AND BAND SERVICES\")\n\nif ($foundStart -and $foundEnd) {\n    $delRange = $d.Range($rngStart.Start, $rngEnd.Start)\n    $delRange.Delete()\n}\n"}
